$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing date-cell formatting (border/font/alignment/number format)
# from the last populated row (A244) down onto the new A245:A247 cells so the
# new rows reuse the same style as the rest of column A.
$ws.Range("A244").Copy()
$ws.Range("A245:A247").PasteSpecial(-4122)

# New data rows to append (column A values are Excel serial date numbers)
$rows = @(
    @{ Row = 245; A = 44319; B = 0; C = 7; D = 132.5757575757576 },
    @{ Row = 246; A = 44320; B = 0; C = 5; D = 94.6969696969697 },
    @{ Row = 247; A = 44321; B = 0; C = 5; D = 94.6969696969697 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
}
